$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the logged hours in D2:D5 (they were entered ~24x too large) ---
$ws.Range("D2").Value = 0.083333333333333329
$ws.Range("D3").Value = 0.083333333333333329
$ws.Range("D4").Value = 0.10416666666666667
$ws.Range("D5").Value = 0.125

# --- Simplify the Totale ore formula in I2 (drop the redundant SUM wrapper) ---
# and give it its own elapsed-time number format.
$ws.Range("I2").Formula = "=F2+G2"
$ws.Range("I2").NumberFormat = "[h]:mm:ss"

# --- D12 no longer needs the underlined font ---
$ws.Range("D12").Font.Underline = $false

# --- New work-log entries dated 15/12/2017 ---
$ws.Range("A14").Value = 43084
$ws.Range("A14").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"
$ws.Range("B14").Value = "Giovanni"
$ws.Range("C14").Value = "Definizione comportamento combinazione, probabilita"
$ws.Range("C14").WrapText = $true
$ws.Range("D14").Value = 0.16666666666666666
$ws.Range("D14").NumberFormat = "h:mm;@"
$ws.Rows("14").RowHeight = 43.75

$ws.Range("A15").Value = 43084
$ws.Range("A15").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"
$ws.Range("B15").Value = "Mirko"
$ws.Range("C15").Value = "Definizione comportamento combinazione, probabilita"
$ws.Range("C15").WrapText = $true
$ws.Range("D15").Value = 0.16666666666666666
$ws.Range("D15").NumberFormat = "h:mm;@"
$ws.Rows("15").RowHeight = 43.75

# --- Widen column I so the "Totale ore" figure fits ---
$ws.Columns("I:I").ColumnWidth = 20.6

# --- Reset the view: scroll back to the top and select the total cell ---
[void]$ws.Range("I2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
